# "Mise à jour de l'application" - enter N3J3 (matchday 3, National 3) playing-time
# data for every player on the roster: minutes played in column DG and the
# T/R/NR/HG (Titulaire / Remplaçant / Non entré / Hors groupe) status in column DH.
# Also corrects Enzo Vita's N3J1 status (CZ4) from "HG" to "NR".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (minutes played or $null if he didn't take the field, status code)
$data = @(
    @{ Row = 2;  Minutes = 90;   Status = "T"  },
    @{ Row = 3;  Minutes = $null; Status = "HG" },
    @{ Row = 4;  Minutes = $null; Status = "NR" },
    @{ Row = 5;  Minutes = $null; Status = "HG" },
    @{ Row = 6;  Minutes = $null; Status = "HG" },
    @{ Row = 7;  Minutes = 59;   Status = "T"  },
    @{ Row = 8;  Minutes = $null; Status = "HG" },
    @{ Row = 9;  Minutes = 90;   Status = "T"  },
    @{ Row = 10; Minutes = $null; Status = "HG" },
    @{ Row = 11; Minutes = 90;   Status = "T"  },
    @{ Row = 12; Minutes = $null; Status = "HG" },
    @{ Row = 13; Minutes = 90;   Status = "T"  },
    @{ Row = 14; Minutes = 75;   Status = "T"  },
    @{ Row = 15; Minutes = 31;   Status = "R"  },
    @{ Row = 16; Minutes = 90;   Status = "T"  },
    @{ Row = 17; Minutes = $null; Status = "HG" },
    @{ Row = 18; Minutes = 59;   Status = "T"  },
    @{ Row = 19; Minutes = 31;   Status = "R"  },
    @{ Row = 20; Minutes = 75;   Status = "T"  },
    @{ Row = 21; Minutes = $null; Status = "HG" },
    @{ Row = 22; Minutes = 90;   Status = "T"  },
    @{ Row = 23; Minutes = $null; Status = "HG" },
    @{ Row = 24; Minutes = 90;   Status = "T"  },
    @{ Row = 25; Minutes = 15;   Status = "R"  },
    @{ Row = 26; Minutes = 15;   Status = "R"  },
    @{ Row = 27; Minutes = $null; Status = "HG" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    if ($null -ne $entry.Minutes) {
        $ws.Range("DG$r").Value = $entry.Minutes
    }
    $ws.Range("DH$r").Value = $entry.Status
}

# Correction: Enzo Vita's N3J1 entry was mis-marked "HG" (hors groupe) and
# should be "NR" (non entré en jeu).
$ws.Range("CZ4").Value = "NR"

Write-Output "N3J3 data entered"
